$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = 45246
$ws.Range("A22").NumberFormat = "d-mmm"
$ws.Range("B22").Value = "Internship"
$ws.Range("C22").Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

$ws.Range("C23").Select() | Out-Null
